$d = $word.ActiveDocument

function Set-ParagraphText($paragraph, $pattern, $replacement) {
    $xml = $paragraph.Range.WordOpenXML
    $xml = $xml -replace $pattern, $replacement
    # InsertXML stamps w:rsidR/w:rsidRDefault on the touched paragraph;
    # strip them back out so the paragraph element matches its original shape.
    $xml = $xml -replace ' w:rsidR="00000000"', ''
    $xml = $xml -replace ' w:rsidRDefault="00000000"', ''
    $paragraph.Range.InsertXML($xml)
}

# 1. "Menu:" -> "Меню:"
Set-ParagraphText $d.Paragraphs.Item(1) '<w:t(?: xml:space="preserve")?>Menu:</w:t>' '<w:t>Меню:</w:t>'

# 2. "Welcome" -> "Добро пожаловать"
Set-ParagraphText $d.Paragraphs.Item(2) '<w:t(?: xml:space="preserve")?>Welcome</w:t>' '<w:t xml:space="preserve">Добро пожаловать</w:t>'

# 3. "What is SmartCash?" -> "Что такое SmartCash?"
Set-ParagraphText $d.Paragraphs.Item(3) '<w:t(?: xml:space="preserve")?>What is SmartCash\?</w:t>' '<w:t xml:space="preserve">Что такое SmartCash?</w:t>'

# 4. "Brochure/Whitepaper" -> "Брошюра"
Set-ParagraphText $d.Paragraphs.Item(4) '<w:t(?: xml:space="preserve")?>Brochure/Whitepaper</w:t>' '<w:t>Брошюра</w:t>'
